$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (and their relationships) before rebuilding the row data
$ws.Hyperlinks.Delete()

# Clear any existing data below the header row so the sheet can be rebuilt cleanly
$ws.Range("A2:H1000").Clear()

# Column widths (this engine stores OOXML width = ColumnWidth + 5/6)
$ws.Columns.Item(1).ColumnWidth = 21 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 51 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 12 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 30 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 12 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 44 - (5/6)
$ws.Columns.Item(7).ColumnWidth = 12 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 21 - (5/6)

# Row 2
$ws.Range("A2").Value = '2025-10-02 12:35:29'
$ws.Range("B2").Value = '【AIで開発生産性を革新】AI活用推進エンジニア募集(副業・業務委託)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5391864'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5391864') | Out-Null
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-10-02 12:35:29'
$ws.Range("B3").Value = 'Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5405426'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5405426') | Out-Null
$ws.Range("G3").Value = 305
$ws.Range("H3").Value = '🔥Python ◆開発,スクレイピング'

# Row 4
$ws.Range("A4").Value = '2025-10-02 12:35:29'
$ws.Range("B4").Value = '【緊急・即日対応歓迎】SIM AI(Docker)のログイン遅延解消とGoogle認証テスト'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5405408'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5405408') | Out-Null
$ws.Range("G4").Value = 298
$ws.Range("H4").Value = '🔥AI,Ai'

# Row 5
$ws.Range("A5").Value = '2025-10-02 12:35:29'
$ws.Range("B5").Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5405023'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5405023') | Out-Null
$ws.Range("G5").Value = 178
$ws.Range("H5").Value = '★bot ◆ツール'

# Row 6
$ws.Range("A6").Value = '2025-10-02 12:35:29'
$ws.Range("B6").Value = '【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5405540'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5405540') | Out-Null
$ws.Range("G6").Value = 158
$ws.Range("H6").Value = '◆自動化,スクレイピング ◇管理'

# Row 7
$ws.Range("A7").Value = '2025-10-02 12:35:29'
$ws.Range("B7").Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5251319') | Out-Null
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = '◆ツール,スクレイピング ◇サイト'

# Row 8
$ws.Range("A8").Value = '2025-10-02 12:35:29'
$ws.Range("B8").Value = '【急募】データ処理のためのExcel VBA・マクロ開発依頼 もしくはスクレイピングによる対応'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5405218'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5405218') | Out-Null
$ws.Range("G8").Value = 98
$ws.Range("H8").Value = '◆開発,スクレイピング'

# Row 9
$ws.Range("A9").Value = '2025-10-02 12:35:29'
$ws.Range("B9").Value = '【急募】集計分析ツール(keyence社製「KI」)の設定構築経験者'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5405052'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5405052') | Out-Null
$ws.Range("G9").Value = 73
$ws.Range("H9").Value = '◆ツール'

# Row 10
$ws.Range("A10").Value = '2025-10-02 12:35:29'
$ws.Range("B10").Value = '【SES案件多数/リモート可】フロントエンドエンジニア募集(HTML/CSS〜モダンFWまで歓迎)'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5399721'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5399721') | Out-Null
$ws.Range("G10").Value = 25

# Row 11
$ws.Range("A11").Value = '2025-10-02 12:35:29'
$ws.Range("B11").Value = '【在宅勤務】ランサーズ業務委託で働ける、ネパール人個人の方を募集します!'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5404906'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5404906') | Out-Null
$ws.Range("G11").Value = 18

# Row 12
$ws.Range("A12").Value = '2025-10-02 12:35:29'
$ws.Range("B12").Value = 'LINE公式アカウントの動作確認・タグ等設定対応'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5405235'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5405235') | Out-Null
$ws.Range("G12").Value = 10

# Restore the Hyperlink cell style on column F (Hyperlinks.Add applies its own xf)
$ws.Range("F2:F12").Style = "Hyperlink"
